# "Generate Report for Archive"
# 1) Status text "Ready for handoff" -> "In Translation" everywhere it appears
#    (Overview!E2:F2/E3:F3 and the Status column (C2:C3) on the zh-cn / de-de sheets).
# 2) Narrow the "Status/zh-cn/de-de" report columns to match the refreshed
#    report layout (was ~17.22 chars, now ~13.41 chars).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E3").Value = "In Translation"
$overview.Range("F3").Value = "In Translation"
$overview.Columns.Item(5).ColumnWidth = 12.42
$overview.Columns.Item(6).ColumnWidth = 12.42

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C3").Value = "In Translation"
$zhcn.Columns.Item(3).ColumnWidth = 12.42

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C3").Value = "In Translation"
$dede.Columns.Item(3).ColumnWidth = 12.42
